$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("commondata")
$ws.Range("B3").Value = "'123456"
